$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 23.11014822181284
$ws.Range("C2").Value = 4.362399951498389
$ws.Range("D2").Value = 4.072610582564346
$ws.Range("E2").Value = 10.98770551718587
$ws.Range("F2").Value = 55.77520609594342
$ws.Range("J2").Value = 10.54774182860877
$ws.Range("K2").Value = 19.0138559828884
$ws.Range("L2").Value = 11.57069076593228
$ws.Range("N2").Value = 26.33251820137313

$ws.Range("B3").Value = 23.01054385223804
$ws.Range("C3").Value = 4.225191698887947
$ws.Range("D3").Value = 4.079697481276648
$ws.Range("E3").Value = 11.00629192790168
$ws.Range("F3").Value = 55.7402908538503
$ws.Range("J3").Value = 10.5652801806097
$ws.Range("K3").Value = 18.94947017374801
$ws.Range("L3").Value = 11.5852058515645
$ws.Range("N3").Value = 26.36910687975005

$ws.Range("B4").Value = 22.95471465579007
$ws.Range("C4").Value = 4.140185687579641
$ws.Range("D4").Value = 4.084579972641775
$ws.Range("E4").Value = 11.01884644128617
$ws.Range("F4").Value = 55.72818307309109
$ws.Range("J4").Value = 10.57678869411716
$ws.Range("K4").Value = 18.91422504216162
$ws.Range("L4").Value = 11.59577601372099
$ws.Range("N4").Value = 26.39343557131997

$ws.Range("B5").Value = 22.93332124749067
$ws.Range("C5").Value = 4.105417302812604
$ws.Range("D5").Value = 4.086703563336673
$ws.Range("E5").Value = 11.02425022335494
$ws.Range("F5").Value = 55.7255980272905
$ws.Range("J5").Value = 10.58166498311184
$ws.Range("K5").Value = 18.90095104056489
$ws.Range("L5").Value = 11.60050059881648
$ws.Range("N5").Value = 26.40381832645444

$ws.Range("B6").Value = 22.92985135312854
$ws.Range("C6").Value = 4.099638247599882
$ws.Range("D6").Value = 4.087064284519724
$ws.Range("E6").Value = 11.02516490779655
$ws.Range("F6").Value = 55.72531071282157
$ws.Range("J6").Value = 10.58248596237768
$ws.Range("K6").Value = 18.89881293580175
$ws.Range("L6").Value = 11.60131031496336
$ws.Range("N6").Value = 26.40557068202578

$ws.Range("B7").Value = 22.95442061871004
$ws.Range("C7").Value = 4.139717219775897
$ws.Range("D7").Value = 4.08460806930987
$ws.Range("E7").Value = 11.01891815302656
$ws.Range("F7").Value = 55.72813869669325
$ws.Range("J7").Value = 10.57685370185924
$ws.Range("K7").Value = 18.91404160353057
$ws.Range("L7").Value = 11.59583804174437
$ws.Range("N7").Value = 26.39357369926628

$ws.Range("B8").Value = 23.07470985126999
$ws.Range("C8").Value = 4.315284923833558
$ws.Range("D8").Value = 4.074944146152868
$ws.Range("E8").Value = 10.99387728485612
$ws.Range("F8").Value = 55.76123281552786
$ws.Range("J8").Value = 10.55363574196723
$ws.Range("K8").Value = 18.99077240310506
$ws.Range("L8").Value = 11.57535166643915
$ws.Range("N8").Value = 26.34474736712688

$ws.Range("B9").Value = 23.3520190219868
$ws.Range("C9").Value = 4.650974767891701
$ws.Range("D9").Value = 4.060190355806073
$ws.Range("E9").Value = 10.95381647650092
$ws.Range("F9").Value = 55.89999844643911
$ws.Range("J9").Value = 10.51395716349023
$ws.Range("K9").Value = 19.17475127010964
$ws.Range("L9").Value = 11.54831855604573
$ws.Range("N9").Value = 26.26377667652612

$ws.Range("B10").Value = 23.57975211656572
$ws.Range("C10").Value = 4.889233984693912
$ws.Range("D10").Value = 4.051886533502564
$ws.Range("E10").Value = 10.92987039635468
$ws.Range("F10").Value = 56.0467014494914
$ws.Range("J10").Value = 10.48834671072442
$ws.Range("K10").Value = 19.32958333394062
$ws.Range("L10").Value = 11.53644792998582
$ws.Range("N10").Value = 26.21328973538323

$ws.Range("B11").Value = 23.68826327015263
$ws.Range("C11").Value = 4.995235995653557
$ws.Range("D11").Value = 4.048654557087459
$ws.Range("E11").Value = 10.92016221349591
$ws.Range("F11").Value = 56.12307334761731
$ws.Range("J11").Value = 10.47745935491517
$ws.Range("K11").Value = 19.40410315043884
$ws.Range("L11").Value = 11.53277721150354
$ws.Range("N11").Value = 26.1922754516654

$ws.Range("B12").Value = 23.73003279284274
$ws.Range("C12").Value = 5.034992114536621
$ws.Range("D12").Value = 4.047508705461141
$ws.Range("E12").Value = 10.91665588911231
$ws.Range("F12").Value = 56.15336953863756
$ws.Range("J12").Value = 10.47344588919453
$ws.Range("K12").Value = 19.43289198901642
$ws.Range("L12").Value = 11.53163527125238
$ws.Range("N12").Value = 26.18459859264247

$ws.Range("B13").Value = 23.72100725126502
$ws.Range("C13").Value = 5.026447645314939
$ws.Range("D13").Value = 4.047752021425952
$ws.Range("E13").Value = 10.91740348731084
$ws.Range("F13").Value = 56.14678369210668
$ws.Range("J13").Value = 10.47430540385769
$ws.Range("K13").Value = 19.42666674826191
$ws.Range("L13").Value = 11.53187018518103
$ws.Range("N13").Value = 26.18623945264625

$ws.Range("B14").Value = 23.69168623614454
$ws.Range("C14").Value = 4.998514686352642
$ws.Range("D14").Value = 4.048558725769997
$ws.Range("E14").Value = 10.91987034255486
$ws.Range("F14").Value = 56.12553832096043
$ws.Range("J14").Value = 10.47712697547642
$ws.Range("K14").Value = 19.40646031093003
$ws.Range("L14").Value = 11.53267829510338
$ws.Range("N14").Value = 26.19163824413329

$ws.Range("B15").Value = 23.6738138044393
$ws.Range("C15").Value = 4.981353696930436
$ws.Range("D15").Value = 4.049063004677143
$ws.Range("E15").Value = 10.92140348214277
$ws.Range("F15").Value = 56.11270377567252
$ws.Range("J15").Value = 10.47886949776773
$ws.Range("K15").Value = 19.39415694632278
$ws.Range("L15").Value = 11.5332055745574
$ws.Range("N15").Value = 26.1949817284641

$ws.Range("B16").Value = 23.57275718296581
$ws.Range("C16").Value = 4.882254870139279
$ws.Range("D16").Value = 4.052108692134213
$ws.Range("E16").Value = 10.93052865615671
$ws.Range("F16").Value = 56.04190345509063
$ws.Range("J16").Value = 10.48907354582315
$ws.Range("K16").Value = 19.3247940129633
$ws.Range("L16").Value = 11.53672256912997
$ws.Range("N16").Value = 26.21470235956148

$ws.Range("B17").Value = 23.51200138025521
$ws.Range("C17").Value = 4.820819368676361
$ws.Range("D17").Value = 4.054116541594458
$ws.Range("E17").Value = 10.93642984729128
$ws.Range("F17").Value = 56.00093116000959
$ws.Range("J17").Value = 10.49552854572483
$ws.Range("K17").Value = 19.28327628831571
$ws.Range("L17").Value = 11.53932264791689
$ws.Range("N17").Value = 26.22730042436364

$ws.Range("B18").Value = 23.47752055678925
$ws.Range("C18").Value = 4.785261180453179
$ws.Range("D18").Value = 4.055322771098443
$ws.Range("E18").Value = 10.93993562243816
$ws.Range("F18").Value = 55.97827263111927
$ws.Range("J18").Value = 10.49931312722738
$ws.Range("K18").Value = 19.25978210564935
$ws.Range("L18").Value = 11.54098095310745
$ws.Range("N18").Value = 26.23473027253612

$ws.Range("B19").Value = 23.46592651885984
$ws.Range("C19").Value = 4.773184987425741
$ws.Range("D19").Value = 4.055740014522412
$ws.Range("E19").Value = 10.94114179278244
$ws.Range("F19").Value = 55.97075701134191
$ws.Range("J19").Value = 10.50060687089892
$ws.Range("K19").Value = 19.25189415068576
$ws.Range("L19").Value = 11.54157040602644
$ws.Range("N19").Value = 26.23727745404485

$ws.Range("B20").Value = 23.51842108379944
$ws.Range("C20").Value = 4.827382587399924
$ws.Range("D20").Value = 4.05389748918328
$ws.Range("E20").Value = 10.93579011198755
$ws.Range("F20").Value = 56.00519885978161
$ws.Range("J20").Value = 10.49483396791953
$ws.Range("K20").Value = 19.28765612658341
$ws.Range("L20").Value = 11.53902901907137
$ws.Range("N20").Value = 26.22594031880381

$ws.Range("B21").Value = 23.70028033803947
$ws.Range("C21").Value = 5.006730021583079
$ws.Range("D21").Value = 4.048319662808194
$ws.Range("E21").Value = 10.91914115849835
$ws.Range("F21").Value = 56.13174134690485
$ws.Range("J21").Value = 10.47629524701413
$ws.Range("K21").Value = 19.4123801133565
$ws.Range("L21").Value = 11.532434205827
$ws.Range("N21").Value = 26.1900448672641

$ws.Range("B22").Value = 23.82307808317368
$ws.Range("C22").Value = 5.121686163110381
$ws.Range("D22").Value = 4.045128840528029
$ws.Range("E22").Value = 10.90925051102123
$ws.Range("F22").Value = 56.22245801346135
$ws.Range("J22").Value = 10.46481626518693
$ws.Range("K22").Value = 19.49720679468905
$ws.Range("L22").Value = 11.52956979682734
$ws.Range("N22").Value = 26.16822181542907

$ws.Range("B23").Value = 23.75718744640092
$ws.Range("C23").Value = 5.060551194476958
$ws.Range("D23").Value = 4.046790383085196
$ws.Range("E23").Value = 10.9144388641244
$ws.Range("F23").Value = 56.17331112558728
$ws.Range("J23").Value = 10.4708846378186
$ws.Range("K23").Value = 19.45163622647935
$ws.Range("L23").Value = 11.53096651889935
$ws.Range("N23").Value = 26.17971941950323

$ws.Range("B24").Value = 23.51551733670434
$ws.Range("C24").Value = 4.824416094629944
$ws.Range("D24").Value = 4.053996361127751
$ws.Range("E24").Value = 10.93607898399508
$ws.Range("F24").Value = 56.00326663791253
$ws.Range("J24").Value = 10.49514775752357
$ws.Range("K24").Value = 19.2856748325651
$ws.Range("L24").Value = 11.53916125939796
$ws.Range("N24").Value = 26.22655463984544

$ws.Range("B25").Value = 23.27268844527934
$ws.Range("C25").Value = 4.561423169918421
$ws.Range("D25").Value = 4.063734655670963
$ws.Range("E25").Value = 10.96368849974102
$ws.Range("F25").Value = 55.85457453530039
$ws.Range("J25").Value = 10.52406749867214
$ws.Range("K25").Value = 19.12146703894784
$ws.Range("L25").Value = 11.55422676585059
$ws.Range("N25").Value = 26.28410005372527

